$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.377.57'
$ws.Range("E2").Value = '  +0.99%  '
$ws.Range("D3").Value = '1.851.06'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  +1.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6194'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.014'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07460'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2959'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.07'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07751'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").Value = '1.831.59'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.025'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6755'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009058'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.910'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.58%  '
$ws.Range("D18").Value = '29.347.13'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = '2.087.01'
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.25%  '
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.016'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.193'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.017'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.50'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1432'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.540'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.509'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.171'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05614'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.124'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.223'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7512'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.855'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("E36").Value = '  +0.98%  '
$ws.Range("E37").Value = '  +2.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.839'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01788'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").Value = '1.216.91'
$ws.Range("E40").Value = '  -1.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.514'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9130'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.015'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.58'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").Value = '1.987.15'
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5160'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.34%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000121'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4067'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.173'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05855'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.95%  '
